$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows (I2:J20)
$data = @(
    @(1,4),
    @(3,6),
    @(5,6),
    @(4,6),
    @(7,8),
    @(9,9),
    @(1,2),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,4),
    @(1,4),
    @(1,4),
    @(1,5),
    @(1,4),
    @(1,4),
    @(1,3),
    @(3,4),
    @(1,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
